$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (COM BGR-packed integers) matching the workbook's existing fill palette
$colorRed    = 255      # FFFF0000
$colorYellow = 65535    # FFFFFF00
$colorGreen1 = 5296274  # FF92D050
$colorGreen2 = 5287936  # FF00B050
$colorOrange = 49407    # FFFFC000 (new fill)

# --- Row 2 ---
$ws.Range("B2").Value = "constructura"
$ws.Range("B2").Interior.Color = $colorRed
$ws.Range("C2").Value = "bono_pie_inmo"
$ws.Range("C2").Interior.Color = $colorYellow

# --- Row 3 ---
$ws.Range("A3").Interior.Color = $colorRed
$ws.Range("B3").Value = "inmobiliaria"
$ws.Range("B3").Interior.Color = $colorGreen1
$ws.Range("C3").Value = "precio_incluye_bono_pie"
$ws.Range("C3").Interior.Color = $colorYellow

# --- Row 4 ---
$ws.Range("A4").Interior.Color = $colorGreen2
$ws.Range("B4").Value = "latitud"
$ws.Range("B4").Interior.Color = $colorGreen1
$ws.Range("C4").Value = "descuento_precio"
$ws.Range("C4").Interior.Color = $colorYellow

# --- Row 5 ---
$ws.Range("B5").Value = "longitud"
$ws.Range("B5").Interior.Color = $colorGreen1

# --- Row 6 ---
$ws.Range("A6").Interior.Color = $colorGreen2
$ws.Range("B6").Value = "descripcion"
$ws.Range("B6").Interior.Color = $colorGreen2

# --- Row 7 ---
$ws.Range("A7").Interior.Color = $colorGreen2

# --- Row 8-13 (left column fills only) ---
$ws.Range("A8").Interior.Color = $colorRed
$ws.Range("A9").Interior.Color = $colorGreen1
$ws.Range("A10").Interior.Color = $colorGreen2
$ws.Range("A11").Interior.Color = $colorGreen1
$ws.Range("A12").Interior.Color = $colorRed
$ws.Range("A13").Interior.Color = $colorGreen1

# --- Row 15: estado_uso -> estado_uso* with orange fill ---
$ws.Range("A15").Value = "estado_uso*"
$ws.Range("A15").Interior.Color = $colorOrange

# --- Row 16 ---
$ws.Range("A16").Value = "fecha_entrega"
$ws.Range("A16").Interior.Color = $colorRed

# --- Row 7 new cells (added after estado_uso* so shared-string order matches) ---
$ws.Range("F7").Value = "*"
$ws.Range("F7").Interior.Color = $colorOrange
$ws.Range("G7").Value = "No especifica si es nuevo pero si tiene el año de construccion"

# --- Selection moves to B2 ---
$ws.Range("B2").Select()
